$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.714.92'
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = '3.333.57'
$ws.Range("E3").Value = '  +1.98%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.23'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.74'
$ws.Range("E6").Value = '  +2.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +2.15%  '
$ws.Range("D9").Value = '3.330.72'
$ws.Range("E9").Value = '  +2.06%  '
$ws.Range("E10").Value = '  +6.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.579'
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.82'
$ws.Range("E12").Value = '  +4.82%  '
$ws.Range("E13").Value = '  +1.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '690.88'
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("D15").Value = '3.877.44'
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("E16").Value = '  +2.57%  '
$ws.Range("D17").Value = '67.679.83'
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").Value = '3.334.06'
$ws.Range("E19").Value = '  +2.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.59'
$ws.Range("E20").Value = '  +2.70%  '
$ws.Range("E21").Value = '  +4.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.894'
$ws.Range("E22").Value = '  +1.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.45'
$ws.Range("E23").Value = '  +4.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.92'
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '101.85'
$ws.Range("E25").Value = '  +4.58%  '
$ws.Range("E26").Value = '  +2.14%  '
$ws.Range("E27").Value = '  +2.01%  '
$ws.Range("E28").Value = '  +6.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.97'
$ws.Range("E29").Value = '  +1.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.54'
$ws.Range("E30").Value = '  +3.65%  '
$ws.Range("E31").Value = '  +6.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '570.87'
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.00'
$ws.Range("E33").Value = '  +1.89%  '
$ws.Range("E34").Value = '  +3.31%  '
$ws.Range("D35").Value = '3.718.10'
$ws.Range("E35").Value = '  -1.90%  '
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '57.18'
$ws.Range("E37").Value = '  +3.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.28'
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.99'
$ws.Range("E39").Value = '  +12.06%  '
$ws.Range("E40").Value = '  +4.78%  '
$ws.Range("E41").Value = '  +7.27%  '
$ws.Range("E42").Value = '  +2.58%  '
$ws.Range("E43").Value = '  -2.18%  '
$ws.Range("E44").Value = '  +3.13%  '
$ws.Range("E45").Value = '  +3.77%  '
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("E47").Value = '  +6.07%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  +0.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '132.05'
$ws.Range("E51").Value = '  +3.68%  '
